$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy style from existing header cell (H1) to the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Data values for column I (I0) and column J (IF)
$i0 = @(5, 4, 5, 5, 1, 1, 8, 8, 6, 1, 6)
$if = @(6, 5, 6, 6, 2, 2, 8, 8, 6, 1, 6)

for ($r = 0; $r -lt $i0.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $i0[$r]
    $ws.Cells.Item($row, 10).Value = $if[$r]
}
